$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new activity ("Inschrijftab implementeren") is being inserted as row 18
# (just above "Testen + Testrapport opstellen"). Insert a fresh blank row
# below the current row 19 so the existing row 19 keeps its original
# (unstyled) formatting and only the new row below picks up row 19's old
# per-cell formatting/content.
$ws.Rows("20:20").Insert()

# Move the old row 19 ("Testen + Testrapport opstellen", ...) down into the
# newly created row 20, preserving its formula (still referencing row 18).
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = $ws.Range("B19").Value2
$ws.Range("C20").Formula = "=(C18+D18)"
$ws.Range("D20").Value = $ws.Range("D19").Value2
$ws.Range("E20").Value = $ws.Range("E19").Value2
$ws.Range("F20").Value = $ws.Range("F19").Value2

# Turn row 19 into the new activity. Only the activity name (B) and actor
# (F) are filled in; the dependency/date/duration/predecessor columns are
# left blank for this entry.
$ws.Range("B19").Value = "Inschrijftab implementeren"
$ws.Range("C19").Clear()
$ws.Range("D19").Clear()
$ws.Range("E19").Clear()
$ws.Range("F19").Value = "Bas"

$ws.Range("F19").Select() | Out-Null
